$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.15000000000049
$ws.Range("H2").Value = [double]"1.714630153861246e-16"
$ws.Range("K2").Value = 56.76729839412141
$ws.Range("L2").Value = "[50.759121408047534, 62.77547538019529]"
$ws.Range("O2").Value = 1.767342413731195
$ws.Range("P2").Value = "[1.6541318676558872, 1.8805529598065034]"
$ws.Range("S2").Value = 57.97555061291111
$ws.Range("T2").Value = "[54.25805646747016, 61.693044758352066]"
$ws.Range("W2").Value = 18.07577577577613
$ws.Range("X2").Value = 17.62262262262297
$ws.Range("Y2").Value = 18.52892892892929

# Row 3 updates
$ws.Range("E3").Value = 22.97000000000015
$ws.Range("H3").Value = [double]"1.714630153861246e-16"
$ws.Range("K3").Value = 51.97659959929638
$ws.Range("L3").Value = "[44.39678402001279, 59.55641517857997]"
$ws.Range("O3").Value = -2.163579324994773
$ws.Range("P3").Value = "[-2.31452671976185, -2.0126319302276956]"
$ws.Range("S3").Value = 57.2821961102841
$ws.Range("T3").Value = "[53.37283282472116, 61.19155939584704]"
$ws.Range("W3").Value = 7.909589589589643
$ws.Range("X3").Value = 7.357757757757808
$ws.Range("Y3").Value = 8.461421421421479
